{"js": "// Office.js (Word JavaScript API) script that applies the\n// \"Added a few more slots\" edit:\n//   1. Insert a new \"Meta description\" paragraph right after the title\n//      (Heading1) paragraph at the top of the document.\n//   2. Remove the duplicate bold title paragraph that was sitting right\n//      before the closing italic paragraph at the end of the document.\n//   3. Replace the text of that trailing italic paragraph with the new\n//      \"Create a feature image ...\" image-prompt copy (keeping the\n//      italic formatting).\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) Insert the \"Meta description: ...\" paragraph after paragraph 1\n//    (the \"Play 3 Stars Slot Game Free - Review & Demo\" Heading1).\n// ---------------------------------------------------------------------\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst titlePara = paras.items[0];\nconst metaText =\n  \"Meta description: Take a spin and play 3 Stars, a 5-reel video slot \" +\n  \"game with 50 fixed paylines and four jackpots available to win. \" +\n  \"Review and demo available.\";\nconst metaPara = titlePara.insertParagraph(metaText, Word.InsertLocation.after);\nmetaPara.style = \"Normal\";\nawait context.sync();\n\n// Bold just the \"Meta description\" label, leaving the\n// \": Take a spin ...\" remainder un-bolded.\nconst labelResults = metaPara.search(\"Meta description\", { matchCase: true });\nlabelResults.load(\"items\");\nawait context.sync();\nlabelResults.items[0].font.bold = true;\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Delete the duplicate bold \"Play 3 Stars Slot Game Free - Review &\n//    Demo\" paragraph near the end of the document. (Match on the\n//    \"Normal\"-styled copy, not the Heading1 title at the top of the\n//    document, which has identical text.)\n// ---------------------------------------------------------------------\nconst allParas = body.paragraphs;\nallParas.load(\"items/text,items/style\");\nawait context.sync();\n\nconst dupeTitle = allParas.items.find(\n  (p) => p.text === \"Play 3 Stars Slot Game Free - Review & Demo\" && p.style === \"Normal\"\n);\nif (dupeTitle) {\n  dupeTitle.delete();\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3) Replace the text of the trailing italic paragraph (previously the\n//    meta-description sentence) with the new image-prompt copy, keeping\n//    its italic formatting intact.\n// ---------------------------------------------------------------------\nconst finalParas = body.paragraphs;\nfinalParas.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph by its (still unique) *exact* old sentence rather\n// than by position, so the script is resilient to any ordering\n// differences. Use an exact match (not a substring match) so the newly\n// inserted \"Meta description: Take a spin ...\" paragraph from step 1 -\n// which contains the same words as a substring - is not picked up by\n// mistake.\nconst oldSentence =\n  \"Take a spin and play 3 Stars, a 5-reel video slot game with 50 fixed \" +\n  \"paylines and four jackpots available to win. Review and demo available.\";\nlet targetPara = finalParas.items.find((p) => p.text === oldSentence);\nif (!targetPara) {\n  targetPara = finalParas.items[finalParas.items.length - 1];\n}\n\nconst newImagePromptText =\n  'Create a feature image for \"3 Stars\" that showcases a happy Maya ' +\n  \"warrior with glasses in a cartoon style. Use bright colors to make \" +\n  \"the image stand out and include elements of Chinese culture to \" +\n  \"depict the theme of the game. The Maya warrior should be holding a \" +\n  \"stack of gold coins, surrounded by Chinese lanterns, and standing \" +\n  'in front of a temple. The title of the game, \"3 Stars\", should be ' +\n  \"prominently displayed in the image, along with the logo of La JVL, \" +\n  \"the game development company.\";\n\nconst targetRange = targetPara.getRange();\ntargetRange.insertText(newImagePromptText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop script (PowerShell-style) that applies the\n# \"Added a few more slots\" edit:\n#   1. Insert a new \"Meta description\" paragraph right after the title\n#      (Heading1) paragraph at the top of the document.\n#   2. Remove the duplicate bold title paragraph that was sitting right\n#      before the closing italic paragraph at the end of the document.\n#   3. Replace the text of that trailing italic paragraph with the new\n#      \"Create a feature image ...\" image-prompt copy (keeping the\n#      italic formatting).\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Insert the \"Meta description: ...\" paragraph after paragraph 1\n#    (the \"Play 3 Stars Slot Game Free - Review & Demo\" Heading1).\n# ---------------------------------------------------------------------\n$titlePara = $d.Paragraphs(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs(2)\n$metaPara.Style = \"Normal\"\n\n$metaRange = $metaPara.Range\n$metaRange.InsertBefore(\"Meta description: Take a spin and play 3 Stars, a 5-reel video slot game with 50 fixed paylines and four jackpots available to win. Review and demo available.\")\n\n# Bold just the \"Meta description\" label (16 characters), leaving the\n# \": Take a spin ...\" remainder un-bolded.\n$labelRange = $d.Range($metaRange.Start, $metaRange.Start + 16)\n$labelRange.Bold = 1\n\n# ---------------------------------------------------------------------\n# 2) Delete the duplicate bold \"Play 3 Stars Slot Game Free - Review &\n#    Demo\" paragraph near the end of the document. (Match on the\n#    \"Normal\"-styled copy, not the Heading1 title at the top of the\n#    document, which has identical text.)\n# ---------------------------------------------------------------------\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -eq \"Play 3 Stars Slot Game Free - Review & Demo`r\" -and $p.Style.NameLocal -eq \"Normal\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# ---------------------------------------------------------------------\n# 3) Replace the text of the trailing italic paragraph (previously the\n#    meta-description sentence) with the new image-prompt copy, keeping\n#    its italic formatting intact. Locate it by its (still unique) exact\n#    old sentence rather than by position, so the script is resilient to\n#    any ordering differences. An exact match (not a substring/-like\n#    match) is required so the newly inserted \"Meta description: Take a\n#    spin ...\" paragraph from step 1 - which contains the same words as\n#    a substring - is not picked up by mistake.\n# ---------------------------------------------------------------------\n$oldSentence = \"Take a spin and play 3 Stars, a 5-reel video slot game with 50 fixed paylines and four jackpots available to win. Review and demo available.`r\"\n$newCount = $d.Paragraphs.Count\n$targetPara = $null\nfor ($i = 1; $i -le $newCount; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -eq $oldSentence) {\n        $targetPara = $p\n        break\n    }\n}\nif ($targetPara -eq $null) {\n    $targetPara = $d.Paragraphs($newCount)\n}\n\n$targetRangeFull = $targetPara.Range\n$textRange = $d.Range($targetRangeFull.Start, $targetRangeFull.End - 1)\n$textRange.Text = \"Create a feature image for \"\"3 Stars\"\" that showcases a happy Maya warrior with glasses in a cartoon style. Use bright colors to make the image stand out and include elements of Chinese culture to depict the theme of the game. The Maya warrior should be holding a stack of gold coins, surrounded by Chinese lanterns, and standing in front of a temple. The title of the game, \"\"3 Stars\"\", should be prominently displayed in the image, along with the logo of La JVL, the game development company.\"\n"}
